$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Quadruple Fanged Panther" -> "Silver Fanged Panther" (row 35)
# and update its description text accordingly. Other cells on the sheet
# keep the same displayed text (only their underlying shared-string
# indices shift once the old strings are removed from the table).
$ws.Range("A35").Value = "Silver Fanged Panther"
$ws.Range("B35").Value = "A rare Silver Fanged Panther that attacks your opponent in battle. A rare drop from Premium Loot Boxes"

# Reflect the final selection left in the sheet view.
$ws.Range("T30").Select()
